# Updated cryptos list on Mon Aug 12 10:41:23 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 26 and 27 swapped coin identity (Binance-Peg BSC-USD <-> WrappedeETH)
# plus new price/volume figures.
$ws.Range("B26").Value = "WrappedeETH"
$ws.Range("C26").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"

$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"

# Price (column D) updates - force Text storage so values like
# "509.79" or "19.00" keep their exact original formatting instead of
# being auto-converted to floating point numbers by Excel.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "58.642.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.605.44"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "509.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.41"
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.629.74"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.42"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.338"
$ws.Range("D12").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.061.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "58.555.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.621.37"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.59"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.40"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("D22").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.66"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.709.88"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.993"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.160"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0823"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.05"
$ws.Range("D30").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.00"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "148.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.863"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.28"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "290.92"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.63"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.42"
$ws.Range("D43").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "19.70"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0540"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.73"
$ws.Range("D49").Style = "Normal"

# Volume(1h) (column E) updates
$ws.Range("E2").Value = "  -3.89%  "
$ws.Range("E3").Value = "  -2.15%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("E5").Value = "  -3.75%  "
$ws.Range("E6").Value = "  -6.00%  "
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("E8").Value = "  -1.24%  "
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("E10").Value = "  -1.18%  "
$ws.Range("E11").Value = "  -4.28%  "
$ws.Range("E12").Value = "  -4.00%  "
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("E14").Value = "  -2.21%  "
$ws.Range("E15").Value = "  -4.01%  "
$ws.Range("E16").Value = "  -4.33%  "
$ws.Range("E17").Value = "  -3.44%  "
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("E19").Value = "  -4.09%  "
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("E21").Value = "  -2.63%  "
$ws.Range("E22").Value = "  -2.49%  "
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("E24").Value = "  -1.67%  "
$ws.Range("E25").Value = "  -2.24%  "
$ws.Range("E26").Value = "  -2.73%  "
$ws.Range("E27").Value = "  -0.61%  "
$ws.Range("E28").Value = "  -4.91%  "
$ws.Range("E29").Value = "  -4.01%  "
$ws.Range("E30").Value = "  -3.81%  "
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("E33").Value = "  -2.73%  "
$ws.Range("E34").Value = "  -4.72%  "
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("E36").Value = "  +11.81%  "
$ws.Range("E37").Value = "  -2.65%  "
$ws.Range("E38").Value = "  -4.73%  "
$ws.Range("E39").Value = "  -6.35%  "
$ws.Range("E40").Value = "  -1.69%  "
$ws.Range("E41").Value = "  -4.83%  "
$ws.Range("E42").Value = "  -3.93%  "
$ws.Range("E43").Value = "  -4.94%  "
$ws.Range("E44").Value = "  -2.01%  "
$ws.Range("E45").Value = "  -4.44%  "
$ws.Range("E47").Value = "  -3.58%  "
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("E49").Value = "  -4.09%  "
$ws.Range("E50").Value = "  -5.36%  "
